# Updated estimation for WalletTransactions user story.
#
# The "GetBalance" user-story block (rows 40-43) becomes the "GetWalletInfo"
# block (the old task texts describing GetBalance work are replaced with the
# matching GetWalletInfo task texts), and a brand-new user story about
# viewing the wallet's transaction list is appended as rows 48-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rework the "GetBalance" tasks into the "GetWalletInfo" tasks ---------
$ws.Range("B40").Value = "Dodavanje GetWalletInfo metode na WalletService"
$ws.Range("B41").Value = "Implementacija testova za GetWalletInfo"
$ws.Range("B42").Value = "Dodavanje rute za dobijanje informacija o Walletu u WalletController"
$ws.Range("B43").Value = "Dodavanje stranice za za dobijanje informacija o Walletu u MVC aplikaciju"

# --- New user story: "Kao korisnik ... vidim listu svojih transakcija" ---
$ws.Range("A48").Value = "Kao korisnik potrebno je da mogu da vidim listu svojih transakcija"

$ws.Range("B49").Value = "Prosirivanje GetWalletInfo sa transakcijama Walleta"
$ws.Range("C49").Value = 5

$ws.Range("B50").Value = "Prosirivanje GetWalletInfo testova sa transakcijama"
$ws.Range("C50").Value = 10

$ws.Range("B51").Value = "Prosirivanje ViewModela za WalletInfo transakcijama I prikaz transakcija na WalletInfo stranici"
$ws.Range("C51").Value = 30

# --- Cosmetic: move the visible window/selection like the author's save ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I56").Select()

Write-Output "WalletTransactions estimation rows updated"
